# Generate Report for Handoff
# Replace the old localization job id / hash tokens and timestamps with the
# newly generated ones across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$oldGuid = "33d7e144-c21a-4153-a64d-c09c7eddf8dc"
$newGuid = "080c98f6-cd1d-4b29-8f16-db87b04b60c9"

$oldZhHash = "255e84e8e93db11dd728702d06415c6d26d92d40"
$newZhHash = "384be33da16516d6d99a43aa507a1e656f45f390"
$oldDeHash = "255e84e8e93db11dd728702d06415c6d26d92d40"
$newDeHash = "384be33da16516d6d99a43aa507a1e656f45f390"

$newFileName = "$newGuid.md"
$newDisplayPath = "e2e\$newGuid.md"

$newZhXlf = "$newGuid.$newZhHash.zh-cn.xlf"
$newDeXlf = "$newGuid.$newDeHash.de-de.xlf"

$newOverviewDate = "2016-08-17 16:56:35"
$newZhDate = "2016-08-17 16:56:30"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("B2").Value = $newDisplayPath
$wsOverview.Hyperlinks.Item(1).TextToDisplay = $newDisplayPath
$wsOverview.Range("G2").Value = $newOverviewDate

# --- zh-cn sheet ---
$wsZhCn.Range("A2").Value = $newFileName
$wsZhCn.Hyperlinks.Item(1).TextToDisplay = $newFileName
$wsZhCn.Range("G2").Value = $newZhXlf
$wsZhCn.Range("H2").Value = $newZhDate

# --- de-de sheet ---
$wsDeDe.Range("A2").Value = $newFileName
$wsDeDe.Hyperlinks.Item(1).TextToDisplay = $newFileName
$wsDeDe.Range("G2").Value = $newDeXlf

Write-Host "Report regenerated for handoff with job id $newGuid"
